$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: grab a Range covering $searchText (first match in the whole doc) to
# use purely as a *formatting template* — we then overwrite its text with the
# text we actually want to insert at the target location while preserving
# the exact run properties (rPr) captured via FormattedText.
# ---------------------------------------------------------------------------
function Get-Template([string]$searchText) {
    $r = $d.Content
    $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    return $r
}

# NOTE: Range objects in this host are NOT "live" — once earlier text in the
# document is inserted/deleted, a previously captured Range's Start/End no
# longer point at the right characters. So templates are re-resolved via a
# fresh Find immediately before every use (cheap — the template strings are
# never themselves edited by this script).
function Insert-Styled([int]$atPos, [string]$text, [string]$templateSearch) {
    $template = Get-Template($templateSearch)
    $tplLen = $template.End - $template.Start
    $ins = $d.Range($atPos, $atPos)
    $ins.FormattedText = $template.FormattedText
    $newRange = $d.Range($atPos, $atPos + $tplLen)
    $newRange.Text = $text
    return $atPos + $text.Length
}

# ===========================================================================
# Location 1 — <head> paragraph: "Boisson <sn>ayant goust de vin</sn></head>"
#            -> "Boisson ayant <sn>goust</sn> de <m>vin</m></head>"
# ===========================================================================

$d.Content.Find.Execute("oisson ", $true, $false, $false, $false, $false, $true, 1, $false, "oisson ayant ", 2) | Out-Null
$d.Content.Find.Execute("ayant goust de vin", $true, $false, $false, $false, $false, $true, 1, $false, "goust", 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute("oisson ayant <sn>goust</sn>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pos = $rng.End
$pos = Insert-Styled $pos " de " "p047v_a1"
$pos = Insert-Styled $pos "<m>" "<sn>"
$pos = Insert-Styled $pos "vin" "p047v_a1"
$pos = Insert-Styled $pos "</m>" "<sn>"

# ===========================================================================
# Location 3 — second "<sn>goust de vin</sn>" (inside second <ab> paragraph)
#            -> "<sn>goust</sn> de <m>vin</m>"
# ===========================================================================

# "goust de vin" is — before editing — a single plain run sandwiched between
# the (already separate) "<sn>" / "</sn>" tag runs, exactly like Location 1,
# so a scoped Find/Replace on just that run's text leaves the tag runs and
# their formatting completely untouched. By this point the only remaining
# "goust de vin" text left in the document is this second occurrence (the
# first one was already collapsed to "goust" above).
$d.Content.Find.Execute("goust de vin", $true, $false, $false, $false, $false, $true, 1, $false, "goust", 2) | Out-Null

$rng2 = $d.Content
$rng2.Find.Execute("il aura <sn>goust</sn>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pos2 = $rng2.End
$pos2 = Insert-Styled $pos2 " de " "p047v_a1"
$pos2 = Insert-Styled $pos2 "<m>" "<sn>"
$pos2 = Insert-Styled $pos2 "vin" "p047v_a1"
$pos2 = Insert-Styled $pos2 "</m>" "<sn>"

Write-Host "done"
